# Apply line-spacing update (1.5 lines / "auto" rule at 18pt -> w:line="360")
# to the first three paragraphs of the document (the "CONG HOA XA HOI CHU
# NGHIA VIET NAM" / "Doc lap - Tu do - Hanh phuc" / title-divider paragraphs).
$d = $word.ActiveDocument

for ($i = 1; $i -le 3; $i++) {
    $p = $d.Paragraphs($i)
    $p.Format.LineSpacingRule = 5   # wdLineSpaceMultiple -> lineRule="auto"
    $p.Format.LineSpacing = 18      # 18pt -> w:line="360"
}
